# Update odds values on Sheet1 (row 3 and row 4) as per the latest
# FlashScore data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3 updates (columns G..X) ---
$ws.Range("G3").Value  = 1.8
$ws.Range("H3").Value  = 3.1
$ws.Range("I3").Value  = 5
$ws.Range("J3").Value  = 2.6
$ws.Range("K3").Value  = 1.95
$ws.Range("L3").Value  = 5.5
$ws.Range("M3").Value  = 1.1
$ws.Range("N3").Value  = 7
$ws.Range("O3").Value  = 1.5
$ws.Range("P3").Value  = 2.5
$ws.Range("Q3").Value  = 1.93
$ws.Range("R3").Value  = 1.93
$ws.Range("S3").Value  = 2.5
$ws.Range("T3").Value  = 1.5
$ws.Range("U3").Value  = 4.1
$ws.Range("V3").Value  = 1.23
$ws.Range("W3").Value  = 5
$ws.Range("X3").Value  = 1.17

# --- Row 3 updates (columns AD..AJ) ---
$ws.Range("AD3").Value = 7
$ws.Range("AF3").Value = 15
$ws.Range("AG3").Value = 19
$ws.Range("AI3").Value = 6.5
$ws.Range("AJ3").Value = 6.5

# --- Row 3 updates (columns AN..AP) ---
$ws.Range("AN3").Value = 10
$ws.Range("AO3").Value = 23
$ws.Range("AP3").Value = 19

# --- Row 4 updates ---
$ws.Range("O4").Value  = 1.29
$ws.Range("P4").Value  = 3.5
$ws.Range("S4").Value  = 1.9
$ws.Range("T4").Value  = 1.95
$ws.Range("U4").Value  = 2.44
$ws.Range("V4").Value  = 1.54

$wb.Save()
